# Add signin and signout components
#
# This script:
#  1. Replaces the short "Detail" text of row 23 (the firebase-authentication
#     service row) with the full implementation notes, and grows that row's
#     height to fit the new multi-line text.
#  2. Appends a new row to the project-history table for "Create Login Page".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Flesh out the Detail cell for the firebase-authentication service ---
$detail = "ng g s ./@core/shared/services/firebaseAuthentication`n" + `
    "implement`n" + `
    "- constructor(public afs: AngularFirestore, public afAuth: AngularFireAuth, public router: Router)`n" + `
    "- canActivate(route: ActivatedRouteSnapshot, state: RouterStateSnapshot): boolean`n" + `
    "- GoogleAuth() {`n" + `
    "- AuthLogin(provider: firebase.default.auth.AuthProvider | GoogleAuthProvider)`n" + `
    "- saveDataWithExpiry(key : string, value: string, ttl : number)`n" + `
    "- getDataWithExpiry(key: string)`n" + `
    "- getFirebaseUser(): FirebaseUser`n" + `
    "- logout() "

$ws.Range("D23").Value = $detail
$ws.Rows.Item(23).RowHeight = 180

# --- 2. Add a new row to the table for the "Create Login Page" step ---
$table = $ws.ListObjects.Item(1)
$newRow = $table.ListRows.Add()

# Match the formatting already used by the other rows that only carry a
# Date / No / Change Title / (blank File Name) combination, e.g. row 2.
$ws.Range("A2:C2").Copy()
$ws.Range("A24:C24").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E24").PasteSpecial(-4122)

$ws.Range("A24").Value = 44986
$ws.Range("B24").Value = 21
$ws.Range("C24").Value = "Create Login Page"

$ws.Range("C24").Select() | Out-Null
